$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "marvel" worksheet after the existing "Reg_page" sheet ---
$newWs = $wb.Worksheets.Add($null, $ws1)
$newWs.Name = "marvel"

# --- Update header row on Reg_page: add a "Status" column in F ---
$ws1.Range("F1").Value = "Status"
$ws1.Range("F1").Interior.ColorIndex = 55

# --- Row 2: Sandesh ---
$ws1.Range("A2").Value = "Sandesh"
$ws1.Range("B2").Value = "matters"
$ws1.Range("C2").Value = "sandesh.matters@gmail.com"
$ws1.Range("D2").Value = "smatter"
$ws1.Range("E2").Value = "sandesh1"
$ws1.Range("F2").Value = "Pass"

# --- Row 3: Suhas ---
$ws1.Range("A3").Value = "Suhas"
$ws1.Range("B3").Value = "matters"
$ws1.Range("C3").Value = "Suhas.matters@gmail.com"
$ws1.Range("C3").NumberFormat = "General"
$ws1.Range("D3").Value = "sumatters"
$ws1.Range("E3").Value = "suhas2"
$ws1.Range("F3").Value = "Pass"

# --- Row 4: Gopi ---
$ws1.Range("A4").Value = "Gopi"
$ws1.Range("B4").Value = "matters"
$ws1.Range("C4").Value = "Gopi.matters@gmail.com"
$ws1.Range("D4").Value = "gmatters"
$ws1.Range("E4").Value = "gopi01"
$ws1.Range("F4").Value = "Pass"

# --- Row 5: numeric sample row ---
$ws1.Range("A5").Value = 1
$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 2
$ws1.Range("E5").Value = 2
$ws1.Range("F5").Value = "Pass"

# --- Autofit the columns that now have new bestfit content ---
$ws1.Columns.Item(1).AutoFit()
$ws1.Columns.Item(6).AutoFit()

# --- Restore Reg_page as the active/selected sheet & selection ---
$ws1.Activate()
$ws1.Range("H10").Select()
